# Weekly update: a new day's price observation is inserted at the top of the
# "Vega Central Mapocho de Santiago - Mandarina" data block (rows 656-753),
# pushing all existing rows in that block down by one row, and the last row
# of the block is carried over into a brand-new last row (754).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 656
$lastRow  = 753

# Shift rows [firstRow .. lastRow] down by one row, working bottom-up so that
# source rows are not clobbered before they are read.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $destRow = $r + 1
    $src = $ws.Range("A" + $r + ":T" + $r)
    $dst = $ws.Range("A" + $destRow + ":T" + $destRow)
    $src.Copy($dst)
}

# Populate the now-vacated first row of the block with the new observation.
$ws.Range("D" + $firstRow).Value = 45015
$ws.Range("L" + $firstRow).Value = "Primera"
$ws.Range("M" + $firstRow).Value = 300
$ws.Range("N" + $firstRow).Value = 10000
$ws.Range("O" + $firstRow).Value = 10000
$ws.Range("P" + $firstRow).Value = 10000
$ws.Range("R" + $firstRow).Value = "Provincia de Limarí"
$ws.Range("S" + $firstRow).Value = 1000
